$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new numeric-looking price needs to stay plain text (matches source inlineStr formatting)
$textCells = @("D8", "D17", "D18", "D19", "D23", "D25", "D26", "D27", "D31", "D37", "D40", "D42", "D43", "D44", "D45", "D48", "D49", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.874.47"
$ws.Range("E2").Value = "  -1.05%  "
$ws.Range("D3").Value = "1.563.11"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "21.71"
$ws.Range("E8").Value = "  -1.59%  "
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("E10").Value = "  -1.32%  "
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").Value = "1.785.59"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "1.571.97"
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("E14").Value = "  -1.05%  "
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("D16").Value = "26.887.15"
$ws.Range("E16").Value = "  -1.04%  "
$ws.Range("D17").Value = "61.27"
$ws.Range("E17").Value = "  -2.87%  "
$ws.Range("D18").Value = "214.73"
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("D19").Value = "7.35"
$ws.Range("E19").Value = "  +1.95%  "
$ws.Range("D20").Value = "0.0₃0679"
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").Value = "9.18"
$ws.Range("E23").Value = "  -2.44%  "
$ws.Range("E24").Value = "  +0.64%  "
$ws.Range("D25").Value = "154.10"
$ws.Range("E25").Value = "  +1.14%  "
$ws.Range("D26").Value = "6.75"
$ws.Range("E26").Value = "  +1.93%  "
$ws.Range("D27").Value = "14.93"
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("E28").Value = "  -0.22%  "
$ws.Range("E29").Value = "  -0.79%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").Value = "1.11"
$ws.Range("E31").Value = "  -3.19%  "
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("D33").Value = "1.403.72"
$ws.Range("E33").Value = "  +2.17%  "
$ws.Range("E34").Value = "  -0.81%  "
$ws.Range("E35").Value = "  -1.64%  "
$ws.Range("E36").Value = "  -1.01%  "
$ws.Range("D37").Value = "0.920"
$ws.Range("E37").Value = "  -2.54%  "
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("D40").Value = "0.811"
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("D42").Value = "0.997"
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("D43").Value = "5.37"
$ws.Range("E43").Value = "  +3.02%  "
$ws.Range("D44").Value = "2.17"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").Value = "63.20"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("E46").Value = "  -2.45%  "
$ws.Range("D47").Value = "1.699.46"
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("D48").Value = "86.32"
$ws.Range("E48").Value = "  +0.99%  "
$ws.Range("D49").Value = "0.0505"
$ws.Range("E49").Value = "  +2.73%  "
$ws.Range("D50").Value = "0.0₇0982"
$ws.Range("E50").Value = "  -1.36%  "
$ws.Range("D51").Value = "0.0946"
$ws.Range("E51").Value = "  +0.60%  "
